$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the trailing rows (7-14) entirely first ---
# This removes the old "FWFT00049 00101" / "B" shared-string refs in those rows
# before we rewrite rows 2-6, and leaves the styled-but-empty E column cells
# (roll-no column keeps its text number format even though blank).
$ws.Range("A7:G14").ClearContents()

# --- Row 2 : lot W2308120178-01 / shade CX / roll 6 (plain number) ---
$ws.Range("A2").Value = "FWFT00039 00001"
$ws.Range("B2").Value = 2308120178
$ws.Range("C2").Value = "W2308120178-01"
$ws.Range("D2").Value = "CX"
$ws.Range("E2").ClearFormats()
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 72.32
$ws.Range("G2").Value = 180

# --- Row 3 : lot W2308120178-01 / shade CX / roll "16" (text) ---
$ws.Range("A3").Value = "FWFT00039 00001"
$ws.Range("B3").Value = 2308120178
$ws.Range("C3").Value = "W2308120178-01"
$ws.Range("D3").Value = "CX"
$ws.Range("E3").Value = "16"
$ws.Range("F3").Value = 65.27
$ws.Range("G3").Value = 180

# --- Row 4 : lot W2306220352-01 / shade CU / roll "1" (text) ---
$ws.Range("A4").Value = "FWFT00039 00001"
$ws.Range("B4").Value = 2306220352
$ws.Range("C4").Value = "W2306220352-01"
$ws.Range("D4").Value = "CU"
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = 35.81
$ws.Range("G4").Value = 167

# --- Row 5 : lot W2306220352-01 / shade CU / roll "2" (text) ---
$ws.Range("A5").Value = "FWFT00039 00001"
$ws.Range("B5").Value = 2306220352
$ws.Range("C5").Value = "W2306220352-01"
$ws.Range("D5").Value = "CU"
$ws.Range("E5").Value = "2"
$ws.Range("F5").Value = 36.24
$ws.Range("G5").Value = 170

# --- Row 6 : lot W2306220352-01 / shade CU / roll "3" (text) ---
$ws.Range("A6").Value = "FWFT00039 00001"
$ws.Range("B6").Value = 2306220352
$ws.Range("C6").Value = "W2306220352-01"
$ws.Range("D6").Value = "CU"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = 37.58
$ws.Range("G6").Value = 172

# --- Update the active selection shown when the sheet is re-opened ---
[void]$ws.Range("C10").Select()
